# Remove the "Income (USD)" column (E) entirely, shifting Age (F->E) and
# Cluster (G->F) one column to the left, then populate the new Cluster
# column with the computed cluster assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("Income (USD)") - this shifts F->E, G->F automatically.
$ws.Range("E1:E11").EntireColumn.Delete()

# Update the Cluster column (now column F) with the new values.
$clusterValues = @(0, 0, 1, 0, 0, 2, 3, 4, 0, 0)
for ($i = 0; $i -lt $clusterValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $clusterValues[$i]
}
